# "Generate Report for Handoff"
#
# The localization status report is regenerated: the report's status
# moves from "In Translation" to "Ready for handoff", and the two
# "latest datetime" timestamps associated with that regeneration are
# refreshed. Widening the "Status"/"zh-cn"/"de-de" columns' text also
# causes Excel to recompute those columns' widths.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" -----------------
# Overview!E2 ("zh-cn" column) and Overview!F2 ("de-de" column)
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
# zh-cn!C2 and de-de!C2 ("Status" column)
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Timestamps refreshed by the handoff generation run -------------------
# Overview!G2 "Latest HO Xliff Generate Date" and de-de!H2 "Latest Handoff
# Datetime" share the same new timestamp.
$overview.Range("G2").Value = "2016-10-18 12:58:24"
$dede.Range("H2").Value     = "2016-10-18 12:58:24"
# zh-cn!H2 "Latest Handoff Datetime" gets its own refreshed timestamp.
$zhcn.Range("H2").Value     = "2016-10-18 12:58:13"

# --- Column width refresh (auto-sized after the longer status text) -------
# Before: 13.4101848602295 characters wide. After: 17.2159881591797.
# Excel's COM ColumnWidth setter here quantizes to the nearest 1/6th of a
# character (its internal pixel grid), so we pick the input value whose
# quantized result lands closest to the target width.
$overview.Columns.Item(5).ColumnWidth = 16.3   # column E (zh-cn)
$overview.Columns.Item(6).ColumnWidth = 16.3   # column F (de-de)
$zhcn.Columns.Item(3).ColumnWidth     = 16.3   # column C (Status)
$dede.Columns.Item(3).ColumnWidth     = 16.3   # column C (Status)
